$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing Program..Updated XML properties content from B:H to C:I
# (copies values + formatting, leaving B:H's original content in place for now).
$ws.Range("B1:H26").Copy($ws.Range("C1:I26"))

# Set the new column B's width (raw stored width 16 == COM width 15 + 1/6).
$ws.Columns.Item(2).ColumnWidth = 15 + 1/6

# New header
$ws.Range("B1").Value = "Namespace"

# Per-row namespace values (method -> namespace), overwriting the stale copy of
# the old "Program" column that Copy left behind in column B.
$namespaces = @(
    "create",    # row 2  create.paragraph
    "create",    # row 3  create.heading
    "create",    # row 4  create.table
    "core",      # row 5  insert
    "lists",     # row 6  lists.insert
    "core",      # row 7  replace
    "core",      # row 8  delete
    "blocks",    # row 9  blocks.delete
    "format",    # row 10 format.apply
    "format",    # row 11 format.fontSize
    "format",    # row 12 format.fontFamily
    "format",    # row 13 format.color
    "format",    # row 14 format.align
    "comments",  # row 15 comments.create
    "comments",  # row 16 comments.patch
    "comments",  # row 17 comments.delete
    "lists",     # row 18 lists.setType
    "lists",     # row 19 lists.indent
    "lists",     # row 20 lists.outdent
    "lists",     # row 21 lists.restart
    "lists",     # row 22 lists.exit
    "tables",    # row 23 tables.insertRow
    "tables",    # row 24 tables.deleteRow
    "tables",    # row 25 tables.insertColumn
    "tables"     # row 26 tables.deleteColumn
)

for ($i = 0; $i -lt $namespaces.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $namespaces[$i]
}
